$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 75,6

$data[0,0] = 'Brian Burns'
$data[0,1] = 'Group1'
$data[0,2] = 2.666666666666667
$data[0,3] = 44.33333333333334
$data[0,4] = 28.33333333333333
$data[0,5] = 16
$data[1,0] = 'Brian Burns'
$data[1,1] = 'Group2'
$data[1,2] = 4.333333333333333
$data[1,3] = 61.33333333333334
$data[1,4] = 36
$data[1,5] = 25.33333333333333
$data[2,0] = 'Brian Burns'
$data[2,1] = 'Difference'
$data[2,2] = 1.666666666666667
$data[2,3] = 17
$data[2,4] = 7.666666666666668
$data[2,5] = 9.333333333333332
$data[3,0] = 'Bud Dupree'
$data[3,1] = 'Group1'
$data[3,2] = 2
$data[3,3] = 38.66666666666666
$data[3,4] = 28.33333333333333
$data[3,5] = 10.33333333333333
$data[4,0] = 'Bud Dupree'
$data[4,1] = 'Group2'
$data[4,2] = 1.666666666666667
$data[4,3] = 28.33333333333333
$data[4,4] = 16.33333333333333
$data[4,5] = 12
$data[5,0] = 'Bud Dupree'
$data[5,1] = 'Difference'
$data[5,2] = -0.3333333333333333
$data[5,3] = -10.33333333333333
$data[5,4] = -12
$data[5,5] = 1.666666666666666
$data[6,0] = 'Chris Board'
$data[6,1] = 'Group1'
$data[6,2] = 0.3333333333333333
$data[6,3] = 27.33333333333333
$data[6,4] = 14.66666666666667
$data[6,5] = 12.66666666666667
$data[7,0] = 'Chris Board'
$data[7,1] = 'Group2'
$data[7,2] = 0.6666666666666666
$data[7,3] = 20.33333333333333
$data[7,4] = 14
$data[7,5] = 6.333333333333333
$data[8,0] = 'Chris Board'
$data[8,1] = 'Difference'
$data[8,2] = 0.3333333333333333
$data[8,3] = -7
$data[8,4] = -0.6666666666666661
$data[8,5] = -6.333333333333333
$data[9,0] = 'Christian Rozeboom'
$data[9,1] = 'Group1'
$data[9,2] = 2.666666666666667
$data[9,3] = 76.66666666666667
$data[9,4] = 45.66666666666666
$data[9,5] = 31
$data[10,0] = 'Christian Rozeboom'
$data[10,1] = 'Group2'
$data[10,2] = 0.6666666666666666
$data[10,3] = 26.66666666666667
$data[10,4] = 14.33333333333333
$data[10,5] = 12.33333333333333
$data[11,0] = 'Christian Rozeboom'
$data[11,1] = 'Difference'
$data[11,2] = -2
$data[11,3] = -50
$data[11,4] = -31.33333333333333
$data[11,5] = -18.66666666666666
$data[12,0] = 'De''Vondre Campbell'
$data[12,1] = 'Group1'
$data[12,2] = 4.333333333333333
$data[12,3] = 124.6666666666667
$data[12,4] = 82
$data[12,5] = 42.66666666666666
$data[13,0] = 'De''Vondre Campbell'
$data[13,1] = 'Group2'
$data[13,2] = 2
$data[13,3] = 83.33333333333333
$data[13,4] = 48.66666666666666
$data[13,5] = 34.66666666666666
$data[14,0] = 'De''Vondre Campbell'
$data[14,1] = 'Difference'
$data[14,2] = -2.333333333333333
$data[14,3] = -41.33333333333334
$data[14,4] = -33.33333333333334
$data[14,5] = -8
$data[15,0] = 'Demario Davis'
$data[15,1] = 'Group1'
$data[15,2] = 8
$data[15,3] = 111.6666666666667
$data[15,4] = 76.66666666666667
$data[15,5] = 35
$data[16,0] = 'Demario Davis'
$data[16,1] = 'Group2'
$data[16,2] = 5.666666666666667
$data[16,3] = 122
$data[16,4] = 65.66666666666667
$data[16,5] = 56.33333333333334
$data[17,0] = 'Demario Davis'
$data[17,1] = 'Difference'
$data[17,2] = -2.333333333333333
$data[17,3] = 10.33333333333333
$data[17,4] = -11
$data[17,5] = 21.33333333333334
$data[18,0] = 'Denzel Perryman'
$data[18,1] = 'Group1'
$data[18,2] = 1.666666666666667
$data[18,3] = 90
$data[18,4] = 59
$data[18,5] = 31
$data[19,0] = 'Denzel Perryman'
$data[19,1] = 'Group2'
$data[19,2] = 1.666666666666667
$data[19,3] = 71.33333333333333
$data[19,4] = 46
$data[19,5] = 25.33333333333333
$data[20,0] = 'Denzel Perryman'
$data[20,1] = 'Difference'
$data[20,2] = 0
$data[20,3] = -18.66666666666667
$data[20,4] = -13
$data[20,5] = -5.666666666666668
$data[21,0] = 'Duke Riley'
$data[21,1] = 'Group1'
$data[21,2] = 0.3333333333333333
$data[21,3] = 30.11111111111111
$data[21,4] = 17.66666666666667
$data[21,5] = 12.44444444444444
$data[22,0] = 'Duke Riley'
$data[22,1] = 'Group2'
$data[22,2] = 1
$data[22,3] = 32.66666666666666
$data[22,4] = 20
$data[22,5] = 12.66666666666667
$data[23,0] = 'Duke Riley'
$data[23,1] = 'Difference'
$data[23,2] = 0.6666666666666667
$data[23,3] = 2.555555555555554
$data[23,4] = 2.333333333333332
$data[23,5] = 0.2222222222222232
$data[24,0] = 'Eric Kendricks'
$data[24,1] = 'Group1'
$data[24,2] = 7.333333333333333
$data[24,3] = 120
$data[24,4] = 73.33333333333333
$data[24,5] = 46.66666666666666
$data[25,0] = 'Eric Kendricks'
$data[25,1] = 'Group2'
$data[25,2] = 5
$data[25,3] = 130.6666666666667
$data[25,4] = 79
$data[25,5] = 51.66666666666666
$data[26,0] = 'Eric Kendricks'
$data[26,1] = 'Difference'
$data[26,2] = -2.333333333333333
$data[26,3] = 10.66666666666666
$data[26,4] = 5.666666666666671
$data[26,5] = 5
$data[27,0] = 'Eric Wilson'
$data[27,1] = 'Group1'
$data[27,2] = 2.888888888888889
$data[27,3] = 71.55555555555556
$data[27,4] = 37.11111111111111
$data[27,5] = 34.44444444444444
$data[28,0] = 'Eric Wilson'
$data[28,1] = 'Group2'
$data[28,2] = 0.6666666666666666
$data[28,3] = 40
$data[28,4] = 21
$data[28,5] = 19
$data[29,0] = 'Eric Wilson'
$data[29,1] = 'Difference'
$data[29,2] = -2.222222222222222
$data[29,3] = -31.55555555555556
$data[29,4] = -16.11111111111111
$data[29,5] = -15.44444444444444
$data[30,0] = 'Frankie Luvu'
$data[30,1] = 'Group1'
$data[30,2] = 0.6666666666666666
$data[30,3] = 26.66666666666667
$data[30,4] = 17.66666666666667
$data[30,5] = 9
$data[31,0] = 'Frankie Luvu'
$data[31,1] = 'Group2'
$data[31,2] = 5.333333333333333
$data[31,3] = 111.6666666666667
$data[31,4] = 65.33333333333333
$data[31,5] = 46.33333333333334
$data[32,0] = 'Frankie Luvu'
$data[32,1] = 'Difference'
$data[32,2] = 4.666666666666666
$data[32,3] = 85
$data[32,4] = 47.66666666666666
$data[32,5] = 37.33333333333334
$data[33,0] = 'Haason Reddick'
$data[33,1] = 'Group1'
$data[33,2] = 3.333333333333333
$data[33,3] = 69
$data[33,4] = 41
$data[33,5] = 28
$data[34,0] = 'Haason Reddick'
$data[34,1] = 'Group2'
$data[34,2] = 1.666666666666667
$data[34,3] = 33.66666666666666
$data[34,4] = 24.33333333333333
$data[34,5] = 9.333333333333334
$data[35,0] = 'Haason Reddick'
$data[35,1] = 'Difference'
$data[35,2] = -1.666666666666667
$data[35,3] = -35.33333333333334
$data[35,4] = -16.66666666666667
$data[35,5] = -18.66666666666666
$data[36,0] = 'Jadeveon Clowney'
$data[36,1] = 'Group1'
$data[36,2] = 3
$data[36,3] = 29
$data[36,4] = 19.66666666666667
$data[36,5] = 9.333333333333334
$data[37,0] = 'Jadeveon Clowney'
$data[37,1] = 'Group2'
$data[37,2] = 4
$data[37,3] = 39
$data[37,4] = 20.66666666666667
$data[37,5] = 18.33333333333333
$data[38,0] = 'Jadeveon Clowney'
$data[38,1] = 'Difference'
$data[38,2] = 1
$data[38,3] = 10
$data[38,4] = 1
$data[38,5] = 8.999999999999998
$data[39,0] = 'jalenreeves-maybin'
$data[39,1] = 'Group1'
$data[39,2] = 1.333333333333333
$data[39,3] = 43
$data[39,4] = 30
$data[39,5] = 13
$data[40,0] = 'jalenreeves-maybin'
$data[40,1] = 'Group2'
$data[40,2] = 0.6666666666666666
$data[40,3] = 16.33333333333333
$data[40,4] = 10.66666666666667
$data[40,5] = 5.666666666666667
$data[41,0] = 'jalenreeves-maybin'
$data[41,1] = 'Difference'
$data[41,2] = -0.6666666666666666
$data[41,3] = -26.66666666666667
$data[41,4] = -19.33333333333334
$data[41,5] = -7.333333333333333
$data[42,0] = 'Jihad Ward'
$data[42,1] = 'Group1'
$data[42,2] = 0.8888888888888888
$data[42,3] = 17.55555555555556
$data[42,4] = 6.333333333333333
$data[42,5] = 11.22222222222222
$data[43,0] = 'Jihad Ward'
$data[43,1] = 'Group2'
$data[43,2] = 1.666666666666667
$data[43,3] = 25.66666666666667
$data[43,4] = 16.33333333333333
$data[43,5] = 9.333333333333334
$data[44,0] = 'Jihad Ward'
$data[44,1] = 'Difference'
$data[44,2] = 0.7777777777777779
$data[44,3] = 8.111111111111111
$data[44,4] = 10
$data[44,5] = -1.888888888888888
$data[45,0] = 'Kamu Grugier-Hill'
$data[45,1] = 'Group1'
$data[45,2] = 1
$data[45,3] = 53
$data[45,4] = 36.33333333333334
$data[45,5] = 16.66666666666667
$data[46,0] = 'Kamu Grugier-Hill'
$data[46,1] = 'Group2'
$data[46,2] = 1.888888888888889
$data[46,3] = 35.33333333333334
$data[46,4] = 22.44444444444444
$data[46,5] = 12.88888888888889
$data[47,0] = 'Kamu Grugier-Hill'
$data[47,1] = 'Difference'
$data[47,2] = 0.8888888888888886
$data[47,3] = -17.66666666666666
$data[47,4] = -13.88888888888889
$data[47,5] = -3.777777777777777
$data[48,0] = 'Luke Rhodes'
$data[48,1] = 'Group1'
$data[48,2] = 0
$data[48,3] = 1.333333333333333
$data[48,4] = 1.333333333333333
$data[48,5] = 0
$data[49,0] = 'Luke Rhodes'
$data[49,1] = 'Group2'
$data[49,2] = 0
$data[49,3] = 2.666666666666667
$data[49,4] = 1.666666666666667
$data[49,5] = 1
$data[50,0] = 'Luke Rhodes'
$data[50,1] = 'Difference'
$data[50,2] = 0
$data[50,3] = 1.333333333333333
$data[50,4] = 0.3333333333333335
$data[50,5] = 1
$data[51,0] = 'Mack Wilson'
$data[51,1] = 'Group1'
$data[51,2] = 3
$data[51,3] = 54.33333333333334
$data[51,4] = 35.66666666666666
$data[51,5] = 18.66666666666667
$data[52,0] = 'Mack Wilson'
$data[52,1] = 'Group2'
$data[52,2] = 3
$data[52,3] = 49.33333333333334
$data[52,4] = 27.66666666666667
$data[52,5] = 21.66666666666667
$data[53,0] = 'Mack Wilson'
$data[53,1] = 'Difference'
$data[53,2] = 0
$data[53,3] = -5
$data[53,4] = -7.999999999999996
$data[53,5] = 3
$data[54,0] = 'Neville Hewitt'
$data[54,1] = 'Group1'
$data[54,2] = 3
$data[54,3] = 89.66666666666667
$data[54,4] = 57.33333333333334
$data[54,5] = 32.33333333333334
$data[55,0] = 'Neville Hewitt'
$data[55,1] = 'Group2'
$data[55,2] = 1
$data[55,3] = 29
$data[55,4] = 16.33333333333333
$data[55,5] = 12.66666666666667
$data[56,0] = 'Neville Hewitt'
$data[56,1] = 'Difference'
$data[56,2] = -2
$data[56,3] = -60.66666666666667
$data[56,4] = -41
$data[56,5] = -19.66666666666667
$data[57,0] = 'Oren Burks'
$data[57,1] = 'Group1'
$data[57,2] = 0
$data[57,3] = 22.66666666666667
$data[57,4] = 15
$data[57,5] = 7.666666666666667
$data[58,0] = 'Oren Burks'
$data[58,1] = 'Group2'
$data[58,2] = 1
$data[58,3] = 41.66666666666666
$data[58,4] = 21.33333333333333
$data[58,5] = 20.33333333333333
$data[59,0] = 'Oren Burks'
$data[59,1] = 'Difference'
$data[59,2] = 1
$data[59,3] = 19
$data[59,4] = 6.333333333333332
$data[59,5] = 12.66666666666666
$data[60,0] = 'Robert Spillane'
$data[60,1] = 'Group1'
$data[60,2] = 1.333333333333333
$data[60,3] = 37.33333333333334
$data[60,4] = 25.33333333333333
$data[60,5] = 12
$data[61,0] = 'Robert Spillane'
$data[61,1] = 'Group2'
$data[61,2] = 5
$data[61,3] = 128.3333333333333
$data[61,4] = 75
$data[61,5] = 53.33333333333334
$data[62,0] = 'Robert Spillane'
$data[62,1] = 'Difference'
$data[62,2] = 3.666666666666667
$data[62,3] = 91
$data[62,4] = 49.66666666666667
$data[62,5] = 41.33333333333334
$data[63,0] = 'T.J. Watt'
$data[63,1] = 'Group1'
$data[63,2] = 7.333333333333333
$data[63,3] = 57.33333333333334
$data[63,4] = 42
$data[63,5] = 15.33333333333333
$data[64,0] = 'T.J. Watt'
$data[64,1] = 'Group2'
$data[64,2] = 5.666666666666667
$data[64,3] = 56
$data[64,4] = 38.33333333333334
$data[64,5] = 17.66666666666667
$data[65,0] = 'T.J. Watt'
$data[65,1] = 'Difference'
$data[65,2] = -1.666666666666666
$data[65,3] = -1.333333333333336
$data[65,4] = -3.666666666666664
$data[65,5] = 2.333333333333334
$data[66,0] = 'Troy Reeder'
$data[66,1] = 'Group1'
$data[66,2] = 2.666666666666667
$data[66,3] = 76.66666666666667
$data[66,4] = 45.66666666666666
$data[66,5] = 31
$data[67,0] = 'Troy Reeder'
$data[67,1] = 'Group2'
$data[67,2] = 0.6666666666666666
$data[67,3] = 26.66666666666667
$data[67,4] = 14.33333333333333
$data[67,5] = 12.33333333333333
$data[68,0] = 'Troy Reeder'
$data[68,1] = 'Difference'
$data[68,2] = -2
$data[68,3] = -50
$data[68,4] = -31.33333333333333
$data[68,5] = -18.66666666666666
$data[69,0] = 'Zach Cunningham'
$data[69,1] = 'Group1'
$data[69,2] = 2
$data[69,3] = 122.6666666666667
$data[69,4] = 79.66666666666667
$data[69,5] = 43
$data[70,0] = 'Zach Cunningham'
$data[70,1] = 'Group2'
$data[70,2] = 1.666666666666667
$data[70,3] = 37.66666666666666
$data[70,4] = 23.66666666666667
$data[70,5] = 14
$data[71,0] = 'Zach Cunningham'
$data[71,1] = 'Difference'
$data[71,2] = -0.3333333333333333
$data[71,3] = -85
$data[71,4] = -56
$data[71,5] = -29
$data[72,0] = 'Zaire Franklin'
$data[72,1] = 'Group1'
$data[72,2] = 0.6666666666666666
$data[72,3] = 21.66666666666667
$data[72,4] = 13
$data[72,5] = 8.666666666666666
$data[73,0] = 'Zaire Franklin'
$data[73,1] = 'Group2'
$data[73,2] = 6
$data[73,3] = 173
$data[73,4] = 100.6666666666667
$data[73,5] = 72.33333333333333
$data[74,0] = 'Zaire Franklin'
$data[74,1] = 'Difference'
$data[74,2] = 5.333333333333333
$data[74,3] = 151.3333333333333
$data[74,4] = 87.66666666666667
$data[74,5] = 63.66666666666666

$ws.Range("A2:F76").Value = $data
